$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 68 (shifts current rows 68-103 down to 69-104,
# carrying their formatting/values along automatically).
$ws.Rows.Item(68).Insert()

# Populate the new row 68 with this week's data (same record as the
# following week, but with the new date and updated prices).
$ws.Range("A68").Value = 1
$ws.Range("B68").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C68").Value = "Arica y Parinacota"
$ws.Range("D68").Value = 44825
$ws.Range("E68").Value = 15
$ws.Range("F68").Value = 100112038
$ws.Range("G68").Value = "Cebollín baby"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 300
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = 1900
$ws.Range("N68").Value = "`$/paquete 1,5 a 2 kilos"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 950
$ws.Range("Q68").Value = 2
$ws.Range("R68").Value = "Hortaliza"
